$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.142.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.260.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.258.40"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.92"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.800.50"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.259.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.195.47"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.83"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.104"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.47"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0722"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.070.68"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "423.27"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.20"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.76"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.76"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.08%  "
